# Update crypto price/volume data (and a few reordered coin rows) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.519.28'
$ws.Range('E2').Value = '  +15.20%  '

$ws.Range('D3').Value = '1.712.11'
$ws.Range('E3').Value = '  +8.13%  '

$ws.Range('D4').Value = "'0.9872"
$ws.Range('E4').Value = '  -1.50%  '

$ws.Range('D5').Value = "'309.54"
$ws.Range('E5').Value = '  +4.42%  '

$ws.Range('D6').Value = "'0.9826"
$ws.Range('E6').Value = '  -1.01%  '

$ws.Range('D7').Value = "'0.3741"
$ws.Range('E7').Value = '  +4.10%  '

$ws.Range('D8').Value = "'50.69"
$ws.Range('E8').Value = '  +23.84%  '

$ws.Range('D9').Value = "'0.3520"
$ws.Range('E9').Value = '  +6.71%  '

$ws.Range('D10').Value = "'1.203"
$ws.Range('E10').Value = '  +8.96%  '

$ws.Range('D11').Value = "'0.07415"
$ws.Range('E11').Value = '  +7.99%  '

$ws.Range('D12').Value = "'0.9830"
$ws.Range('E12').Value = '  -1.68%  '

$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').Value = "'21.29"
$ws.Range('E13').Value = '  +11.68%  '

$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = "'6.277"
$ws.Range('E14').Value = '  +9.20%  '

$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = "'6.857"
$ws.Range('E15').Value = '  +6.51%  '

$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '1.704.23'
$ws.Range('E16').Value = '  +7.52%  '

$ws.Range('D17').Value = "'0.00001132"
$ws.Range('E17').Value = '  +7.59%  '

$ws.Range('D18').Value = "'0.9826"
$ws.Range('E18').Value = '  -0.99%  '

$ws.Range('D19').Value = "'0.06707"
$ws.Range('E19').Value = '  +2.53%  '

$ws.Range('D20').Value = "'84.92"
$ws.Range('E20').Value = '  +12.80%  '

$ws.Range('D21').Value = "'16.92"
$ws.Range('E21').Value = '  +7.71%  '

$ws.Range('D22').Value = "'6.231"
$ws.Range('E22').Value = '  +6.63%  '

$ws.Range('D23').Value = "'12.38"
$ws.Range('E23').Value = '  +9.00%  '

$ws.Range('D24').Value = '25.380.35'
$ws.Range('E24').Value = '  +14.71%  '

$ws.Range('D25').Value = "'2.390"
$ws.Range('E25').Value = '  +1.28%  '

$ws.Range('D26').Value = "'2.786"
$ws.Range('E26').Value = '  +12.57%  '

$ws.Range('D27').Value = "'153.46"
$ws.Range('E27').Value = '  +3.64%  '

$ws.Range('D28').Value = "'20.14"
$ws.Range('E28').Value = '  +6.62%  '

$ws.Range('D29').Value = '1.888.93'
$ws.Range('E29').Value = '  +7.61%  '

$ws.Range('D30').Value = "'130.17"
$ws.Range('E30').Value = '  +7.24%  '

$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = "'1.158"
$ws.Range('E31').Value = '  +28.13%  '

$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = "'6.519"
$ws.Range('E32').Value = '  +13.20%  '

$ws.Range('B33').Value = 'HuobiToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D33').Value = "'4.002"
$ws.Range('E33').Value = '  +1.74%  '

$ws.Range('B34').Value = 'WEMIXTOKEN'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').Value = "'1.755"
$ws.Range('E34').Value = '  +9.34%  '

$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').Value = "'0.08488"
$ws.Range('E35').Value = '  +5.21%  '

$ws.Range('D36').Value = "'13.13"
$ws.Range('E36').Value = '  +13.77%  '

$ws.Range('D37').Value = "'0.06557"
$ws.Range('E37').Value = '  +10.17%  '

$ws.Range('B38').Value = 'FraxShare'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D38').Value = "'9.076"
$ws.Range('E38').Value = '  +9.94%  '

$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').Value = "'5.447"
$ws.Range('E39').Value = '  +8.24%  '

$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = "'0.02369"
$ws.Range('E40').Value = '  +10.12%  '

$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').Value = "'0.2149"
$ws.Range('E41').Value = '  +9.93%  '

$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = "'1.255"
$ws.Range('E42').Value = '  +2.88%  '

$ws.Range('D43').Value = "'0.6334"
$ws.Range('E43').Value = '  +11.66%  '

$ws.Range('D44').Value = "'0.9823"
$ws.Range('E44').Value = '  -1.00%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = "'13.47"
$ws.Range('E45').Value = '  +7.97%  '

$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = "'0.6107"
$ws.Range('E46').Value = '  +10.86%  '

$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').Value = "'3.812"
$ws.Range('E47').Value = '  +2.09%  '

$ws.Range('D48').Value = "'129.58"
$ws.Range('E48').Value = '  +5.22%  '

$ws.Range('D49').Value = "'2.070"
$ws.Range('E49').Value = '  +8.19%  '

$ws.Range('D50').Value = "'0.07419"
$ws.Range('E50').Value = '  +10.76%  '

$ws.Range('D51').Value = "'77.50"
$ws.Range('E51').Value = '  +8.41%  '
